$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '30.418.75'
Set-TextValue "E2" '  +0.56%  '
Set-TextValue "D3" '1.866.94'
Set-TextValue "E3" '  +0.03%  '
Set-TextValue "D4" '1.001'
Set-TextValue "E4" '  +0.07%  '
Set-TextValue "D5" '247.10'
Set-TextValue "E5" '  +1.54%  '
Set-TextValue "D6" '1.001'
Set-TextValue "E6" '  +0.05%  '
Set-TextValue "D7" '0.4728'
Set-TextValue "E7" '  +0.22%  '
Set-TextValue "D8" '0.2908'
Set-TextValue "E8" '  +1.47%  '
Set-TextValue "D9" '0.06468'
Set-TextValue "E9" '  +0.09%  '
Set-TextValue "D10" '22.00'
Set-TextValue "E10" '  +5.51%  '
Set-TextValue "D11" '0.07710'
Set-TextValue "E11" '  -0.22%  '
Set-TextValue "D12" '97.17'
Set-TextValue "E12" '  +2.38%  '
Set-TextValue "D13" '0.7393'
Set-TextValue "E13" '  +4.44%  '
Set-TextValue "D14" '1.868.87'
Set-TextValue "E14" '  +0.12%  '
Set-TextValue "D15" '5.135'
Set-TextValue "E15" '  +0.93%  '
Set-TextValue "D16" '272.68'
Set-TextValue "E16" '  +0.71%  '
Set-TextValue "D17" '30.427.24'
Set-TextValue "E17" '  +0.64%  '
Set-TextValue "D18" '13.36'
Set-TextValue "E18" '  +0.26%  '
Set-TextValue "D19" '0.9999'
Set-TextValue "E19" '  -0.04%  '
Set-TextValue "D20" '0.000007495'
Set-TextValue "E20" '  -0.51%  '
Set-TextValue "D21" '2.117.80'
Set-TextValue "E21" '  +0.28%  '
Set-TextValue "D22" '1.001'
Set-TextValue "E22" '  +0.05%  '
Set-TextValue "D23" '5.228'
Set-TextValue "E23" '  +0.52%  '
Set-TextValue "D24" '6.149'
Set-TextValue "E24" '  +0.44%  '
Set-TextValue "D25" '9.252'
Set-TextValue "E25" '  -0.77%  '
Set-TextValue "D26" '163.16'
Set-TextValue "E26" '  -1.16%  '
Set-TextValue "D27" '18.68'
Set-TextValue "E27" '  -0.80%  '
Set-TextValue "D28" '1.913'
Set-TextValue "E28" '  +0.37%  '
Set-TextValue "D29" '0.09981'
Set-TextValue "E29" '  +1.35%  '
Set-TextValue "D30" '1.367'
Set-TextValue "E30" '  -0.73%  '
Set-TextValue "E31" '  -0.03%  '
Set-TextValue "D32" '4.256'
Set-TextValue "E32" '  +0.42%  '
Set-TextValue "D33" '4.102'
Set-TextValue "E33" '  +2.18%  '
Set-TextValue "D34" '0.04791'
Set-TextValue "E34" '  +1.30%  '
Set-TextValue "D35" '1.113'
Set-TextValue "E35" '  -0.53%  '
Set-TextValue "D36" '0.6916'
Set-TextValue "E36" '  +0.16%  '
Set-TextValue "D37" '2.711'
Set-TextValue "D38" '0.01846'
Set-TextValue "E38" '  +0.20%  '
Set-TextValue "E39" '  +0.10%  '
Set-TextValue "D40" '6.220'
Set-TextValue "E40" '  -1.41%  '
Set-TextValue "D41" '72.51'
Set-TextValue "E41" '  +3.49%  '
Set-TextValue "D42" '1.967'
Set-TextValue "E42" '  +4.01%  '
Set-TextValue "D43" '0.4172'
Set-TextValue "E43" '  +2.43%  '
Set-TextValue "E44" '  +0.07%  '
Set-TextValue "D45" '0.8344'
Set-TextValue "E45" '  -0.52%  '
Set-TextValue "D46" '101.92'
Set-TextValue "E46" '  -0.05%  '
Set-TextValue "D47" '9.286'
Set-TextValue "E47" '  +1.03%  '
Set-TextValue "D48" '35.38'
Set-TextValue "E48" '  +1.82%  '
Set-TextValue "D49" '6.928'
Set-TextValue "E49" '  -1.86%  '
Set-TextValue "D50" '917.14'
Set-TextValue "E50" '  -0.71%  '
Set-TextValue "D51" '0.05635'
Set-TextValue "E51" '  +1.19%  '
